# Add a new row (32) to the Leetcode summary sheet for
# "138. Copy List with Random Pointer", following the same layout/format
# as the existing Linked List rows (30 and 31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (fill/border/wrap) of the row above so the new row
# matches the sheet's established look (A=plain, B=yellow highlight, C=plain).
$ws.Range("A31:C31").Copy()
$ws.Range("A32:C32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A32").Value2 = "Linked List"
$ws.Range("B32").Value2 = "138. Copy List with Random Pointer"
$ws.Range("C32").Value2 = "This is VERY TRICKY Problem, watch this video: https://youtu.be/5Y2EiZST97Y?si=Hxezj_ZugvAcGz3O`nThe trick is to use a HashMap<Node,Node> we store mapping of oldNode, newNode in it. `nfirst oldNode = head then use while(oldNode!=null) to iter over original LL, make new copy of each node with same val, leave the .random ptr as null. While creating new copies in the loop do hashmap.put(oldNode, newNode).... be sure to store new head of copy LL too`nNow we need to do a 2nd pass over original LL oldNode = head, using while(oldNode!=null) do newNode = hashmap.get(oldNode); newNode.random = hm.get(oldNode.random) || null; `nUSING HASHMAP FOR MAPPING OLDNODE: NEWNODE IS EXTREMLY IMPORTANT!!!!!!!"

# Match the row height Excel computed for the wrapped text (158.4pt == 11 lines).
$ws.Rows.Item(32).RowHeight = 158.4

# Reflect where the author ended up after typing the new description.
$ws.Range("B32").Select() | Out-Null
